$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values stay numeric, just change the numbers
$ws.Range("A2").Value = 4000
$ws.Range("B2").Value = 4000
$ws.Range("C2").Value = 20

# Row 3 values become text (same style as before - no special formatting)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1000.0"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1000.0"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "40"
$ws.Range("C3").Style = "Normal"

# Rows 4 through 6 are removed entirely
$ws.Range("A4:C6").EntireRow.Delete()
